$d = $word.ActiveDocument

# Originally "<id>p080r_N</id>" (N=1,2,3) was split across three separate
# runs:
#   run1: "<id>"     (Courier New, color 7f6000, size 9)
#   run2: "p080r_N"  (default body formatting, color 000000)
#   run3: "</id>"    (Courier New, color 7f6000, size 9)
# The author merged these into a single run containing the full text
# "<id>p080r_N</id>", using run1's (and run3's) formatting. The
# "<id>fig_p080r_N</id>" paragraphs elsewhere in the document are left
# untouched, so we must only target the bare "p080r_N" ids.

for ($i = 1; $i -le 3; $i++) {
    $needle = "<id>p080r_$i</id>"

    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

    if ($found) {
        $matchStart = $rng.Start
        $matchEnd = $rng.End

        # "<id>" is always 4 characters and keeps its original formatting;
        # collapse the remaining "p080r_N</id>" into that first run by
        # deleting it and re-inserting it right after "<id>".
        $firstRun = $d.Range($matchStart, $matchStart + 4)
        $rest = $d.Range($matchStart + 4, $matchEnd)
        $restText = $rest.Text
        $rest.Delete()
        $firstRun.InsertAfter($restText)
    }
}
